$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the XPath string to a path-based comparison
$ws.Range("P5").Value = "//*[local-name()='CelsiusToFahrenheitResult']/text()"

# Remove the explicit (bold/navy) styling from O4 and P5 -> back to default style
$ws.Range("O4").ClearFormats()
$ws.Range("P5").ClearFormats()

# Set width for column P (16th column) so the stored OOXML <col width="..."/>
# comes out to 44.0 (the runtime's ColumnWidth setter rounds to whole pixels
# and re-adds the 5px padding, i.e. stored = (round(chars*6)+5)/6, so an input
# of 43.1667 lands exactly on stored width 44).
$ws.Columns.Item(16).ColumnWidth = 43.1667
